$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data block (header-styled "a1"/"b2" rows with many numeric
# columns) occupies rows 7:8 across columns A:FI. Duplicate that block
# into the next two rows (9:10), matching the original values/styles.
$src = $ws.Range("A7:FI8")
$dst = $ws.Range("A9:FI10")
$src.Copy($dst)
